$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 2019 column header (AE1), matching style/formatting of existing header row (AD1):
# bold + centered, and stored as text (years are text labels, not numbers).
$ws.Range("AE1").NumberFormat = "@"
$ws.Range("AE1").Font.Bold = $true
$ws.Range("AE1").HorizontalAlignment = -4108
$ws.Range("AE1").Value = "2019"

# Add 2019 data values for each sector row
$ws.Range("AE2").Value = 824
$ws.Range("AE3").Value = 166
$ws.Range("AE4").Value = 18800
$ws.Range("AE5").Value = 461
$ws.Range("AE6").Value = 105
$ws.Range("AE7").Value = 923
$ws.Range("AE8").Value = 3776
$ws.Range("AE9").Value = 1667
$ws.Range("AE10").Value = 3732
$ws.Range("AE11").Value = 1674
